# "Generate Report for Handback" — localization-status.xlsx
#
# Flips the Overview/zh-cn/de-de status from "Ready for handoff" to
# "Handed back: in sync with en-US", records the handback xliff files and
# timestamps for each locale, adds a hyperlink to the handback markdown on
# each locale sheet, and widens a few report columns so the new, longer
# values aren't truncated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handbackMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f072ba4d050f6748566b0840897a547c3b2f2cf9/e2e/db9892b6-94ff-490e-84fe-bd693381ec92.md"
$handbackMdName = "db9892b6-94ff-490e-84fe-bd693381ec92.md"

# --- Status columns: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn locale sheet: Latest Target File / Latest Handback File / DateTime ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $handbackMdUrl, $null, $null, $handbackMdName) | Out-Null
$zhcn.Range("J2").Value = "db9892b6-94ff-490e-84fe-bd693381ec92.32c31f72454ee2cae101c4252c2ca760bdd8e5c7.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 00:56:28"

# --- de-de locale sheet: Latest Target File / Latest Handback File / DateTime ---
$dede.Hyperlinks.Add($dede.Range("I2"), $handbackMdUrl, $null, $null, $handbackMdName) | Out-Null
$dede.Range("J2").Value = "db9892b6-94ff-490e-84fe-bd693381ec92.32c31f72454ee2cae101c4252c2ca760bdd8e5c7.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 00:56:35"

# --- Widen columns to fit the longer handback status/file values ---
# ColumnWidth is in characters; Excel snaps to whole-pixel widths the same
# way the real UI does, so these inputs are chosen to land on the closest
# achievable pixel width to the target column widths.
$overview.Columns.Item(5).ColumnWidth = 29.1   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.1   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth = 29.1    # C: Status
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664   # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664  # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.1    # C: Status
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664   # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664  # J: Latest Handback File
